$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 302; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}

$cols = @("S","T","V","W","X","Y")
for ($r = 2; $r -le 24; $r++) {
    $aVal = $ws.Range("A$r").Text
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$r")
        $f = $cell.Formula
        $newF = $f.Substring(0, $f.Length - 1) + ', "' + $aVal + '")'
        $cell.Formula = $newF
    }
}
